$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number, date range) ---
$volCell = $ws.Range("A8")
$volChars = $volCell.Characters(21, 2)
$volChars.Text = "34"

$dateCell = $ws.Range("C9")
$dateChars1 = $dateCell.Characters(27, 9)
$dateChars1.Text = "8/19/2024"
$dateChars2 = $dateCell.Characters(47, 9)
$dateChars2.Text = "8/25/2024"

# --- Cells that switch FROM a numeric style TO the text placeholder style (s=14) ---
# string 20 = "0" placeholder, string 21 = "***.*" placeholder
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("D25"))
$ws.Range("E14").Copy($ws.Range("E25"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("C30"))

# --- Cells that switch FROM the text placeholder style TO a numeric style ---
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("L14").Copy($ws.Range("E20"))

# set final numeric values for the cells that changed style
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100

# --- Plain value updates (style unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("L15").Value = 7.692307692307
$ws.Range("M15").Value = 27.272727272727
$ws.Range("N15").Value = -48.148148148148
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 133.333333333333
$ws.Range("F16").Value = 16
$ws.Range("H16").Value = 77.777777777777
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 20
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = 5
$ws.Range("N16").Value = -79.710144927536
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 6.896551724137
$ws.Range("I17").Value = 223
$ws.Range("J17").Value = 163
$ws.Range("K17").Value = 36.809815950920
$ws.Range("L17").Value = 42.038216560509
$ws.Range("M17").Value = 114.423076923077
$ws.Range("N17").Value = -33.630952380952
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 49
$ws.Range("L18").Value = 11.363636363636
$ws.Range("M18").Value = -45.555555555555
$ws.Range("N18").Value = -90.648854961832
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 45.454545454545
$ws.Range("I19").Value = 108
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = -10.743801652892
$ws.Range("L19").Value = -14.285714285714
$ws.Range("M19").Value = 77.049180327868
$ws.Range("N19").Value = -49.767441860465
$ws.Range("I20").Value = 48
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 9.090909090909
$ws.Range("L20").Value = 26.315789473684
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -84.158415841584
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 84.615384615384
$ws.Range("G21").Value = 59
$ws.Range("H21").Value = 37.288135593220
$ws.Range("I21").Value = 528
$ws.Range("J21").Value = 456
$ws.Range("K21").Value = 15.789473684210
$ws.Range("L21").Value = 20
$ws.Range("M21").Value = 28.467153284671
$ws.Range("N21").Value = -71.17903930131
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -27.272727272727
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 57
$ws.Range("J23").Value = 59
$ws.Range("K23").Value = -3.389830508474
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 119.230769230769
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 57.142857142857
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 19.047619047619
$ws.Range("I24").Value = 408
$ws.Range("J24").Value = 373
$ws.Range("K24").Value = 9.383378016085
$ws.Range("L24").Value = 4.884318766066
$ws.Range("M24").Value = 60.629921259842
$ws.Range("F25").Value = 4
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = -33.333333333333
$ws.Range("L25").Value = -50
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = -6.976744186046
$ws.Range("I26").Value = 280
$ws.Range("J26").Value = 283
$ws.Range("K26").Value = -1.060070671378
$ws.Range("L26").Value = 8.949416342412
$ws.Range("M26").Value = 4.089219330855
$ws.Range("G27").Value = 1
$ws.Range("L27").Value = -4.166666666666
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 28
$ws.Range("K28").Value = -3.448275862068
$ws.Range("L28").Value = -17.647058823529
$ws.Range("L29").Value = -54.166666666666
$ws.Range("N29").Value = -66.666666666666
$ws.Range("L30").Value = -50
$ws.Range("N30").Value = -70.967741935483
